$wb = $excel.ActiveWorkbook

# --- Update selection on the existing "procedimientos" sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("C11:C12").Select()

# --- Add the new "Objetos" sheet after the existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "Objetos"

# --- Header row (bold) ---
$ws2.Cells.Item(1,1).Value = "Procedimiento"
$ws2.Cells.Item(1,2).Value = "Schema"
$ws2.Cells.Item(1,3).Value = "Tablas_involucradas"

# --- Data rows ---
$ws2.Cells.Item(2,1).Value  = "proceso.genera_rdc01"
$ws2.Cells.Item(2,2).Value  = "reporte"
$ws2.Cells.Item(2,3).Value  = "rdc01_texto"

$ws2.Cells.Item(3,1).Value  = "proceso.genera_rdc01"
$ws2.Cells.Item(3,2).Value  = "reporte"
$ws2.Cells.Item(3,3).Value  = "rdc01_detalle"

$ws2.Cells.Item(4,1).Value  = "proceso.genera_rdc01"
$ws2.Cells.Item(4,2).Value  = "reporte"
$ws2.Cells.Item(4,3).Value  = "rdc01_final"

$ws2.Cells.Item(5,1).Value  = "proceso.genera_rdc01"
$ws2.Cells.Item(5,2).Value  = "interface"
$ws2.Cells.Item(5,3).Value  = "cartera_operaciones"

$ws2.Cells.Item(6,1).Value  = "proceso.genera_rdc01"
$ws2.Cells.Item(6,2).Value  = "interno"
$ws2.Cells.Item(6,3).Value  = "tipo_persona_rel"

$ws2.Cells.Item(7,1).Value  = "proceso.genera_rdc01"
$ws2.Cells.Item(7,2).Value  = "interno"
$ws2.Cells.Item(7,3).Value  = "operacion_titulo_rel"

$ws2.Cells.Item(8,1).Value  = "proceso.genera_rdc01"
$ws2.Cells.Item(8,2).Value  = "interno"
$ws2.Cells.Item(8,3).Value  = "tabla_banco_126_rel"

$ws2.Cells.Item(9,1).Value  = "proceso.genera_rdc01"
$ws2.Cells.Item(9,2).Value  = "interface"
$ws2.Cells.Item(9,3).Value  = "tipo_cambio"

$ws2.Cells.Item(10,1).Value = "proceso.genera_rdc01"
$ws2.Cells.Item(10,2).Value = "interface"
$ws2.Cells.Item(10,3).Value = "cuadro_operaciones"

$ws2.Cells.Item(11,1).Value = "proceso.genera_rdc01"
$ws2.Cells.Item(11,2).Value = "interface"
$ws2.Cells.Item(11,3).Value = "cartera_garantias"

$ws2.Cells.Item(12,1).Value = "proceso.genera_rdc01"
$ws2.Cells.Item(12,3).Value = "log_eventos"
$ws2.Cells.Item(12,2).Value = "log"

$ws2.Cells.Item(13,1).Value = "proceso.genera_rdc01"
$ws2.Cells.Item(13,2).Value = "interno"
$ws2.Cells.Item(13,3).Value = "parametros_generales"

$ws2.Cells.Item(14,1).Value = "proceso.genera_rdc01"
$ws2.Cells.Item(14,2).Value = "reporte"
$ws2.Cells.Item(14,3).Value = "rdc01_hist"

# --- Bold the header row ---
$ws2.Range("A1:C1").Font.Bold = $true

# --- Autofit columns A and C like the source workbook ---
$ws2.Columns.Item(1).AutoFit() | Out-Null
$ws2.Columns.Item(3).AutoFit() | Out-Null

# --- Selection on the new sheet (also makes it the active/selected tab) ---
$ws2.Range("E19").Select()
